# Apply cryptos list price/volume updates (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-number-looking price cells to stay text (avoids Excel's
# automatic number conversion truncating/changing the displayed text).
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '59.778.36'
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('D3').Value = '2.670.54'
$ws.Range('E3').Value = '  +2.85%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '538.49'
$ws.Range('E5').Value = '  +0.78%  '
$ws.Range('D6').Value = '145.68'
$ws.Range('E6').Value = '  +3.84%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('E8').Value = '  +1.20%  '
$ws.Range('D9').Value = '2.670.03'
$ws.Range('E9').Value = '  +2.38%  '
$ws.Range('D10').Value = '6.67'
$ws.Range('E10').Value = '  +3.16%  '
$ws.Range('E11').Value = '  +1.04%  '
$ws.Range('E12').Value = '  +1.02%  '
$ws.Range('E13').Value = '  -0.63%  '
$ws.Range('D14').Value = '3.132.98'
$ws.Range('E14').Value = '  +2.52%  '
$ws.Range('D15').Value = '59.730.97'
$ws.Range('E15').Value = '  +0.97%  '
$ws.Range('D16').Value = '21.22'
$ws.Range('E16').Value = '  +3.56%  '
$ws.Range('D17').Value = '2.641.01'
$ws.Range('E17').Value = '  +1.63%  '
$ws.Range('E18').Value = '  +1.53%  '
$ws.Range('D19').Value = '344.99'
$ws.Range('E19').Value = '  -0.17%  '
$ws.Range('E20').Value = '  +2.33%  '
$ws.Range('D21').Value = '10.38'
$ws.Range('E21').Value = '  +2.73%  '
$ws.Range('D22').Value = '6.34'
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').Value = '66.64'
$ws.Range('E24').Value = '  -0.74%  '
$ws.Range('E25').Value = '  +2.56%  '
$ws.Range('E26').Value = '  -1.04%  '
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('D28').Value = '7.32'
$ws.Range('E28').Value = '  +1.83%  '
$ws.Range('D29').Value = '0.0₃0754'
$ws.Range('E29').Value = '  +2.82%  '
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('E31').Value = '  +2.02%  '
$ws.Range('D32').Value = '5.85'
$ws.Range('E32').Value = '  +0.80%  '
$ws.Range('D33').Value = '19.03'
$ws.Range('E33').Value = '  +1.24%  '
$ws.Range('D34').Value = '150.34'
$ws.Range('E34').Value = '  +0.59%  '
$ws.Range('E35').Value = '  +1.68%  '
$ws.Range('E36').Value = '  +2.76%  '
$ws.Range('D37').Value = '0.844'
$ws.Range('E37').Value = '  +1.09%  '
$ws.Range('E38').Value = '  -0.20%  '
$ws.Range('D39').Value = '0.830'
$ws.Range('E39').Value = '  +0.89%  '
$ws.Range('D40').Value = '293.56'
$ws.Range('E40').Value = '  +6.18%  '
$ws.Range('E41').Value = '  +2.21%  '
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.19%  '
$ws.Range('D43').Value = '0.608'
$ws.Range('E43').Value = '  +1.64%  '
$ws.Range('D44').Value = '19.62'
$ws.Range('E44').Value = '  +6.05%  '
$ws.Range('E45').Value = '  +4.27%  '
$ws.Range('D46').Value = '10.75'
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('E47').Value = '  -1.22%  '
$ws.Range('D48').Value = '1.983.92'
$ws.Range('E48').Value = '  +2.33%  '
$ws.Range('E49').Value = '  +1.80%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').Value = '4.57'
$ws.Range('E50').Value = '  +1.67%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = '18.47'
$ws.Range('E51').Value = '  +0.94%  '
